$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column D (new quarterly data), shifting old D:K to G:N
$ws.Columns("D:F").Insert()

# Copy cell formatting (number format, font, alignment) from column G (the first of the
# shifted-original columns) into the newly inserted D:F columns so the new quarter data
# matches the existing look (date format on the header row, #,##0 on data rows, etc.)
$ws.Range("G7:G102").Copy()
$ws.Range("D7:F102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new/changed cell values (columns D through N) for every affected row.
$rowData = @{
  7 = @{ D=43519; E=43428; F=43337; G=43246; H=43155; I=43064; J=42973; K=42882; L=42791; M=42700; N=42609 }
  8 = @{ D=432700; E=493600; F=536200; G=562300; H=468400; I=450000; J=454900; K=476400; L=370500; M=245300; N=263300 }
  9 = @{ D=366300; E=422700; F=452400; G=476700; H=400700; I=387200; J=381400; K=405600; L=321200; M=216400; N=231400 }
  10 = @{ D=66400; E=70900; F=83800; G=85600; H=67700; I=62800; J=73500; K=70800; L=49300; M=28900; N=31900 }
  11 = @{  }
  12 = @{ D="NA"; E="NA"; F="NA"; G="NA"; H="NA"; I="NA"; J="NA"; K="NA"; L="NA"; M="NA"; N="NA" }
  13 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  14 = @{ D=200; E="NA"; F="NA"; G="NA"; H="NA"; I="NA"; J=0; K="NA"; L="NA"; M="NA"; N=0 }
  15 = @{ D=2300; E=2700; F=3400; G=1900; H=1900; I=2100; J=2100; K=10200; L=10400; M=2100; N="NA" }
  16 = @{  }
  17 = @{ D=403800; E=461000; F=490500; G=514000; H=433100; I=418800; J=411500; K=441500; L=342100; M=226900; N=244400 }
  18 = @{ D=28900; E=32600; F=45700; G=48300; H=35300; I=31200; J=43400; K=34900; L=28400; M=18400; N=18900 }
  19 = @{  }
  20 = @{ D=200; E=800; F=300; G=100; H=-100; I=100; J=200; K=0; L=0; M=100; N=200 }
  21 = @{ D="NA"; E=39200; F=52500; G=52700; H=39400; I=35500; J=47800; K=46900; L=40600; M=22100; N=20700 }
  22 = @{ D=4300; E=4500; F=4400; G=4200; H=4900; I=4800; J=5300; K=5300; L=5200; M=1100; N="NA" }
  23 = @{ D=24800; E=28900; F=41600; G=44200; H=30300; I=26500; J=38400; K=29600; L=23200; M=17400; N=19100 }
  24 = @{ D=3200; E=6700; F=11800; G=11700; H=8200; I=8600; J=13500; K=10300; L=7900; M=5600; N=6000 }
  25 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  26 = @{ D=21600; E=22200; F=29800; G=32500; H=22100; I=18000; J=24900; K=19400; L=15300; M=11700; N=13100 }
  27 = @{ D=21600; E=22200; F=29800; G=32500; H=22100; I=18000; J=24900; K=19400; L=15300; M=11700; N=13100 }
  28 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  29 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K="NA"; L="NA"; M="NA"; N="NA" }
  30 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  31 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  32 = @{ D=-200; E=-800; F=-300; G=-100; H=100; I=-100; J=-200; K=0; L=0; M=-100; N=-200 }
  33 = @{ D=21600; E=22200; F=29800; G=32500; H=22100; I=18000; J=24900; K=19400; L=15300; M=11700; N=13100 }
  34 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  35 = @{ D=21600; E=22200; F=29800; G=32500; H=22100; I=18000; J=24900; K=19400; L=15300; M=11700; N=13100 }
  38 = @{ D=43519; E=43428; F=43337; G=43246; H=43155; I=43064; J=42973; K=42882; L=42791; M=42700; N=42609 }
  39 = @{  }
  40 = @{  }
  41 = @{ D=3000; E=700; F=2300; G=39000; H=27400; I=54500; J=35900; K=24400; L=10900; M=25600; N=85600 }
  42 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  43 = @{ D=179900; E=140800; F=164600; G=148900; H=157400; I=116800; J=124500; K=121000; L=120500; M=81800; N=66200 }
  44 = @{ D=189600; E=191500; F=195100; G=177400; H=178000; I=151800; J=142300; K=144400; L=148500; M=155400; N=122500 }
  45 = @{ D=15200; E=10300; F=9900; G=8400; H=9800; I=12400; J=11400; K=8500; L=13900; M=10600; N=6300 }
  46 = @{ D=387800; E=343300; F=371900; G=373800; H=372700; I=335400; J=314100; K=298300; L=293800; M=273400; N=280600 }
  47 = @{ D=27000; E=26700; F=28300; G=28100; H=27900; I=27600; J=27400; K=27000; L=26900; M=26700; N=26500 }
  48 = @{ D=119800; E=112200; F=103200; G=84500; H=80800; I=76100; J=73600; K=70700; L=69900; M=68700; N=55900 }
  49 = @{ D=533900; E=536100; F=538100; G=465200; H=467100; I=469100; J=469200; K=473900; L=484100; M=500300; N=1200 }
  50 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  51 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  52 = @{ D=9800; E=11700; F=10300; G=14100; H=16800; I=16600; J=18200; K=20500; L=20100; M=13700; N=26500 }
  53 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  54 = @{ D=1078100; E=1030000; F=1051800; G=965700; H=965300; I=924900; J=902500; K=890400; L=894700; M=882700; N=390700 }
  55 = @{  }
  56 = @{  }
  57 = @{ D=92600; E=79700; F=81000; G=88400; H=99700; I=76100; J=79200; K=79600; L=66900; M=51000; N=44100 }
  58 = @{ D=2800; E="NA"; F=0; G=0; H=0; I=2300; J=2900; K=12100; L=11300; M=7600; N="NA" }
  59 = @{ D=117200; E=119500; F=123100; G=102000; H=95900; I=100100; J=85100; K=85800; L=73500; M=69900; N=48800 }
  60 = @{ D=212500; E=199200; F=204200; G=190400; H=195600; I=178500; J=167200; K=177500; L=151700; M=128500; N=92900 }
  61 = @{ D=274200; E=253300; F=291400; G=251800; H=271100; I=268400; J=271700; K=274800; L=318200; M=334700; N=0 }
  62 = @{ D=20800; E=21000; F=21700; G=17700; H=20800; I=21300; J=22000; K=21800; L=22300; M=22000; N=29400 }
  63 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  64 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  65 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  66 = @{ D=507500; E=473500; F=517400; G=459900; H=487500; I=468200; J=460800; K=474100; L=492100; M=485200; N=122400 }
  67 = @{  }
  68 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  69 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  70 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  71 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  72 = @{ D=805900; E=787800; F=768800; G=742100; H=712800; I=693900; J=679100; K=654200; L=641200; M=629100; N=620500 }
  73 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  74 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  75 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  76 = @{ D=570700; E=556500; F=534400; G=505800; H=477800; I=456600; J=441700; K=416300; L=402600; M=397500; N=268400 }
  77 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  80 = @{ D=43519; E=43428; F=43337; G=43246; H=43155; I=43064; J=42973; K=42882; L=42791; M=42700; N=42609 }
  81 = @{ D=21600; E=22200; F=29800; G=32500; H=22100; I=18000; J=24900; K=19400; L=15300; M=11700; N=13100 }
  82 = @{  }
  83 = @{ D="NA"; E=5800; F=6600; G=4300; H=4100; I=4200; J=4100; K=12000; L=12200; M=3600; N=1500 }
  84 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  85 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  86 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  87 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  88 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  89 = @{ D="NA"; E=54200; F=22300; G=46000; H=-14400; I=29500; J=29800; K=62200; L=5300; M=-200; N=21000 }
  90 = @{  }
  91 = @{ D="NA"; E=-12800; F=-10500; G=-6400; H=-6300; I=-5400; J=-4300; K=-2800; L=-3400; M=-3600; N=-4600 }
  92 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  93 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  94 = @{ D="NA"; E=-13200; F=-93900; G=-6500; H=-6100; I=-5300; J=-1900; K=-2600; L=-3500; M=-397500; N=-3900 }
  95 = @{  }
  96 = @{ D=0; E=-3200; F=-3200; G=-3200; H=-6400; I=0; J=-3200; K=-3200; L=-3200; M=-3200; N=-2700 }
  97 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  98 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  99 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  100 = @{ D="NA"; E=-42600; F=34900; G=-27900; H=-6500; I=-5600; J=-16300; K=-46200; L=-16500; M=337600; N=-2800 }
  101 = @{ D=0; E=0; F=0; G=0; H=0; I=0; J=0; K=0; L=0; M=0; N=0 }
  102 = @{ D="NA"; E=-1600; F=-36700; G=11600; H=-27000; I=18500; J=11600; K=13400; L=-14700; M=-60000; N=14300 }
}

foreach ($r in $rowData.Keys) {
  $cols = $rowData[$r]
  foreach ($c in $cols.Keys) {
    $ws.Range("$c$r").Value2 = $cols[$c]
  }
}

Write-Host "Done. Dimension:" $ws.UsedRange.Address()